$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 6209.06
$ws.Range("D2").Value2 = 2438.59
$ws.Range("E2").Value2 = 9849.469999999999
$ws.Range("F2").Value2 = 2841.77
$ws.Range("G2").Value2 = 6516.9

$ws.Range("C3").Value2 = 2504.19
$ws.Range("D3").Value2 = 4586.76
$ws.Range("E3").Value2 = 16411.4
$ws.Range("F3").Value2 = 19733.52
$ws.Range("G3").Value2 = 7978.06

$ws.Range("C4").Value2 = 12745.35
$ws.Range("D4").Value2 = 13015.69
$ws.Range("E4").Value2 = 9987.969999999999
$ws.Range("F4").Value2 = 11039.67
$ws.Range("G4").Value2 = 5804.47

$ws.Range("C5").Value2 = 15607.12
$ws.Range("D5").Value2 = 8748.85
$ws.Range("E5").Value2 = 2380.09
$ws.Range("F5").Value2 = 4211.45
$ws.Range("G5").Value2 = 2268.14

$ws.Range("C6").Value2 = 16458.58
$ws.Range("D6").Value2 = 11384.59
$ws.Range("E6").Value2 = 5389.63
$ws.Range("F6").Value2 = 8523.16
$ws.Range("G6").Value2 = 12050.34

$ws.Range("C7").Value2 = 15126.52
$ws.Range("D7").Value2 = 14483.7
$ws.Range("E7").Value2 = 1211.32
$ws.Range("F7").Value2 = 10350.09
$ws.Range("G7").Value2 = 6815.98

$ws.Range("C8").Value2 = 16250.27
$ws.Range("D8").Value2 = 16889.78
$ws.Range("E8").Value2 = 7582.78
$ws.Range("F8").Value2 = 14188.34
$ws.Range("G8").Value2 = 12550.75

$ws.Range("C9").Value2 = 13118.25
$ws.Range("D9").Value2 = 3236.57
$ws.Range("E9").Value2 = 11064.06
$ws.Range("F9").Value2 = 12368.96
$ws.Range("G9").Value2 = 4472.27

$ws.Range("C10").Value2 = 3039.06
$ws.Range("D10").Value2 = 7638.24
$ws.Range("E10").Value2 = 1852.09
$ws.Range("F10").Value2 = 17182.61
$ws.Range("G10").Value2 = 13935.04

$ws.Range("C11").Value2 = 6412.3
$ws.Range("D11").Value2 = 19625.47
$ws.Range("E11").Value2 = 19487.57
$ws.Range("F11").Value2 = 12342.25
$ws.Range("G11").Value2 = 4454.72

$ws.Range("C12").Value2 = 3535.74
$ws.Range("D12").Value2 = 8490.379999999999
$ws.Range("E12").Value2 = 19449.86
$ws.Range("F12").Value2 = 4560.64
$ws.Range("G12").Value2 = 4706.57

$ws.Range("C13").Value2 = 12142.71
$ws.Range("D13").Value2 = 7903.93
$ws.Range("E13").Value2 = 4614.75
$ws.Range("F13").Value2 = 8070.16
$ws.Range("G13").Value2 = 9924.120000000001

$ws.Range("C14").Value2 = 15372.36
$ws.Range("D14").Value2 = 1854.75
$ws.Range("E14").Value2 = 6641.87
$ws.Range("F14").Value2 = 12389.89
$ws.Range("G14").Value2 = 18153.84

$ws.Range("C15").Value2 = 14677.83
$ws.Range("D15").Value2 = 6592.65
$ws.Range("E15").Value2 = 15656.12
$ws.Range("F15").Value2 = 15161.9
$ws.Range("G15").Value2 = 15741.98

$ws.Range("C16").Value2 = 4309.8
$ws.Range("D16").Value2 = 16909.42
$ws.Range("E16").Value2 = 14349.65
$ws.Range("F16").Value2 = 10959.23
$ws.Range("G16").Value2 = 9862.02

$ws.Range("C17").Value2 = 1139.04
$ws.Range("D17").Value2 = 4575.12
$ws.Range("E17").Value2 = 6805.18
$ws.Range("F17").Value2 = 12852.32
$ws.Range("G17").Value2 = 13204.45

$ws.Range("C18").Value2 = 6593.82
$ws.Range("D18").Value2 = 9062.540000000001
$ws.Range("E18").Value2 = 6990.37
$ws.Range("F18").Value2 = 18110.54
$ws.Range("G18").Value2 = 8957.57

$ws.Range("C19").Value2 = 13933.94
$ws.Range("D19").Value2 = 11170.61
$ws.Range("E19").Value2 = 9733.35
$ws.Range("F19").Value2 = 1260.66
$ws.Range("G19").Value2 = 18264.95

$ws.Range("C20").Value2 = 6092.97
$ws.Range("D20").Value2 = 5681
$ws.Range("E20").Value2 = 9607.17
$ws.Range("F20").Value2 = 12513.76
$ws.Range("G20").Value2 = 7977.11

$ws.Range("C21").Value2 = 16730.66
$ws.Range("D21").Value2 = 14466.34
$ws.Range("E21").Value2 = 16337.14
$ws.Range("F21").Value2 = 19388.2
$ws.Range("G21").Value2 = 13109.39

$ws.Range("C22").Value2 = 5050.11
$ws.Range("D22").Value2 = 4322.81
$ws.Range("E22").Value2 = 4593.02
$ws.Range("F22").Value2 = 19693.29
$ws.Range("G22").Value2 = 9708.139999999999

$ws.Range("C23").Value2 = 17866
$ws.Range("D23").Value2 = 8086.23
$ws.Range("E23").Value2 = 16660.29
$ws.Range("F23").Value2 = 18532.17
$ws.Range("G23").Value2 = 2644.8

$ws.Range("C24").Value2 = 17479.13
$ws.Range("D24").Value2 = 3219.05
$ws.Range("E24").Value2 = 4919.78
$ws.Range("F24").Value2 = 8307.16
$ws.Range("G24").Value2 = 16943.83

$ws.Range("C25").Value2 = 15924.83
$ws.Range("D25").Value2 = 12921.57
$ws.Range("E25").Value2 = 6894.23
$ws.Range("F25").Value2 = 7928.03
$ws.Range("G25").Value2 = 2523.25

$ws.Range("C26").Value2 = 19029.97
$ws.Range("D26").Value2 = 17394.73
$ws.Range("E26").Value2 = 7126.07
$ws.Range("F26").Value2 = 18658.57
$ws.Range("G26").Value2 = 7411.15

$ws.Range("C27").Value2 = 2259.36
$ws.Range("D27").Value2 = 11798.18
$ws.Range("E27").Value2 = 14712.95
$ws.Range("F27").Value2 = 12789.74
$ws.Range("G27").Value2 = 15877.06

$ws.Range("C28").Value2 = 4349.92
$ws.Range("D28").Value2 = 15633.46
$ws.Range("E28").Value2 = 5032.5
$ws.Range("F28").Value2 = 9247.82
$ws.Range("G28").Value2 = 9297.26

$ws.Range("C29").Value2 = 17027.45
$ws.Range("D29").Value2 = 7204.4
$ws.Range("E29").Value2 = 19631.37
$ws.Range("F29").Value2 = 15271.96
$ws.Range("G29").Value2 = 17212.11

$ws.Range("C30").Value2 = 4461.11
$ws.Range("D30").Value2 = 5677.79
$ws.Range("E30").Value2 = 7606.77
$ws.Range("F30").Value2 = 19156.11
$ws.Range("G30").Value2 = 14032.13

$ws.Range("C31").Value2 = 1484.51
$ws.Range("D31").Value2 = 2793.68
$ws.Range("E31").Value2 = 16918.19
$ws.Range("F31").Value2 = 14216.09
$ws.Range("G31").Value2 = 16251.54

$ws.Range("C32").Value2 = 6568.33
$ws.Range("D32").Value2 = 18452.06
$ws.Range("E32").Value2 = 18687.21
$ws.Range("F32").Value2 = 13811.38
$ws.Range("G32").Value2 = 19742.49

$ws.Range("C33").Value2 = 9532.450000000001
$ws.Range("D33").Value2 = 16906.1
$ws.Range("E33").Value2 = 18553.21
$ws.Range("F33").Value2 = 19236.62
$ws.Range("G33").Value2 = 15110.19

$ws.Range("C34").Value2 = 11676.17
$ws.Range("D34").Value2 = 6759.29
$ws.Range("E34").Value2 = 13928.27
$ws.Range("F34").Value2 = 8863.84
$ws.Range("G34").Value2 = 1606.57

$ws.Range("C35").Value2 = 7925.52
$ws.Range("D35").Value2 = 1253
$ws.Range("E35").Value2 = 18441.99
$ws.Range("F35").Value2 = 12006.5
$ws.Range("G35").Value2 = 9368.59

$ws.Range("C36").Value2 = 8505.35
$ws.Range("D36").Value2 = 17822.97
$ws.Range("E36").Value2 = 19539.84
$ws.Range("F36").Value2 = 12787.56
$ws.Range("G36").Value2 = 3122.55

